# Apply weekly fruit/vegetable price update: rows 2-26 on Sheet1 have their
# Fecha/Calidad/Volumen/Precio.../Unidad columns (D,I,J,K,L,M,N,P,Q) redistributed
# across rows (row 12 is unchanged). Values below are the ORIGINAL (pre-edit)
# workbook contents, captured once so the permutation reads consistent source data
# even though the loop below writes into the very same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orig = @{}
$orig[2] = @{ D=44648; I='Primera'; J=120; K=6500; L=7000; M=6750; N='$/caja 60 unidades'; P=112; Q=60 }
$orig[3] = @{ D=44967; I='Segunda'; J=50; K=4500; L=5000; M=4850; N='$/caja 90 unidades'; P=54; Q=90 }
$orig[4] = @{ D=44935; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
$orig[5] = @{ D=44242; I='Primera'; J=160; K=5000; L=5500; M=5250; N='$/caja 60 unidades'; P=88; Q=60 }
$orig[6] = @{ D=44421; I='Primera'; J=100; K=8000; L=9000; M=8500; N='$/caja 60 unidades'; P=142; Q=60 }
$orig[7] = @{ D=44627; I='Primera'; J=120; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
$orig[8] = @{ D=44740; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
$orig[9] = @{ D=44760; I='Primera'; J=130; K=7000; L=7500; M=7250; N='$/caja 60 unidades'; P=121; Q=60 }
$orig[10] = @{ D=44281; I='Primera'; J=120; K=5500; L=6000; M=5750; N='$/caja 60 unidades'; P=96; Q=60 }
$orig[11] = @{ D=44963; I='Primera'; J=130; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
$orig[12] = @{ D=44676; I='Primera'; J=120; K=4000; L=4500; M=4250; N='$/caja 60 unidades'; P=71; Q=60 }
$orig[13] = @{ D=44400; I='Primera'; J=120; K=9000; L=10000; M=9500; N='$/caja 60 unidades'; P=158; Q=60 }
$orig[14] = @{ D=44827; I='Primera'; J=120; K=6000; L=7000; M=6500; N='$/caja 60 unidades'; P=108; Q=60 }
$orig[15] = @{ D=45079; I='Primera'; J=130; K=4000; L=5000; M=4462; N='$/caja 60 unidades'; P=74; Q=60 }
$orig[16] = @{ D=44764; I='Primera'; J=120; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; P=125; Q=60 }
$orig[17] = @{ D=44382; I='Primera'; J=160; K=7000; L=8000; M=7438; N='$/caja 60 unidades'; P=124; Q=60 }
$orig[18] = @{ D=45243; I='Primera'; J=120; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; P=125; Q=60 }
$orig[19] = @{ D=44785; I='Primera'; J=130; K=7000; L=8000; M=7500; N='$/caja 60 unidades'; P=125; Q=60 }
$orig[20] = @{ D=44669; I='Primera'; J=130; K=4500; L=5000; M=4750; N='$/caja 60 unidades'; P=79; Q=60 }
$orig[21] = @{ D=44494; I='Primera'; J=120; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; P=92; Q=60 }
$orig[22] = @{ D=45044; I='Primera'; J=190; K=4000; L=5000; M=4526; N='$/caja 60 unidades'; P=75; Q=60 }
$orig[23] = @{ D=44657; I='Primera'; J=100; K=5000; L=5500; M=5250; N='$/caja 60 unidades'; P=88; Q=60 }
$orig[24] = @{ D=44589; I='Primera'; J=110; K=5000; L=6000; M=5500; N='$/caja 60 unidades'; P=92; Q=60 }
$orig[25] = @{ D=44603; I='Primera'; J=140; K=5500; L=6000; M=5750; N='$/caja 60 unidades'; P=96; Q=60 }
$orig[26] = @{ D=44362; I='Primera'; J=120; K=8000; L=9000; M=8500; N='$/caja 60 unidades'; P=142; Q=60 }

# Destination row -> source row (source = where the data now at "row" used to live)
$map = @{}
$map[2] = 4
$map[3] = 11
$map[4] = 24
$map[5] = 13
$map[6] = 3
$map[7] = 17
$map[8] = 6
$map[9] = 21
$map[10] = 18
$map[11] = 10
$map[12] = 12
$map[13] = 8
$map[14] = 23
$map[15] = 9
$map[16] = 2
$map[17] = 20
$map[18] = 22
$map[19] = 16
$map[20] = 26
$map[21] = 5
$map[22] = 25
$map[23] = 19
$map[24] = 14
$map[25] = 15
$map[26] = 7

foreach ($row in 2..26) {
    $src = $map[$row]
    $block = $orig[$src]
    $ws.Cells.Item($row, 4).Value = $block.D    # D: Fecha
    $ws.Cells.Item($row, 9).Value = $block.I    # I: Calidad
    $ws.Cells.Item($row, 10).Value = $block.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $block.K   # K: Precio mínimo
    $ws.Cells.Item($row, 12).Value = $block.L   # L: Precio máximo
    $ws.Cells.Item($row, 13).Value = $block.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $block.N   # N: Unidad de comercialización
    $ws.Cells.Item($row, 16).Value = $block.P   # P: Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $block.Q   # Q: Kg o Unidades
}
